# Update "Muc do hoan thanh" (completion level) values in the grading sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PhieuVanDap")

$ws.Range("D13").Value = 1
$ws.Range("D43").Value = 1
$ws.Range("D44").Value = 1
$ws.Range("D45").Value = 0.5
$ws.Range("D47").Value = 1
